$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated LR-pair results (Cadm3-Cadm1) with an added "M2" cluster
# and refreshed ligand/receptor expression statistics, per Dr Hou advice.

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cadm3"
$ws.Range("C2").Value = "Cadm1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.199582666666667
$ws.Range("H2").Value = 12.598748
$ws.Range("I2").Value = 0.3470882463138872
$ws.Range("J2").Value = 0.3745187327618799
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.515984
$ws.Range("N2").Value = 7.547952
$ws.Range("O2").Value = 0.08781336966822693
$ws.Range("P2").Value = 0.09884082726736673
$ws.Range("Q2").Value = 10.56608279601067
$ws.Range("R2").Value = 95.094745164096
$ws.Range("S2").Value = 0.03047898848105798
$ws.Range("T2").Value = 0.03701774137331005

# Row 3: ECs -> M1
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cadm3"
$ws.Range("C3").Value = "Cadm1"
$ws.Range("D3").Value = "M1"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.199582666666667
$ws.Range("H3").Value = 12.598748
$ws.Range("I3").Value = 0.3470882463138872
$ws.Range("J3").Value = 0.3745187327618799
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 10.10016866666667
$ws.Range("N3").Value = 30.300506
$ws.Range("O3").Value = 0.3525180783492434
$ws.Range("P3").Value = 0.3967867150797739
$ws.Range("Q3").Value = 42.41649326294311
$ws.Range("R3").Value = 381.748439366488
$ws.Range("S3").Value = 0.1223548816081804
$ws.Range("T3").Value = 0.148604057708426

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cadm3"
$ws.Range("C4").Value = "Cadm1"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.199582666666667
$ws.Range("H4").Value = 12.598748
$ws.Range("I4").Value = 0.3470882463138872
$ws.Range("J4").Value = 0.3745187327618799
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.445583666666667
$ws.Range("N4").Value = 19.336751
$ws.Range("O4").Value = 0.2249650320703493
$ws.Range("P4").Value = 0.2532157683969216
$ws.Range("Q4").Value = 27.06876144308312
$ws.Range("R4").Value = 243.618852987748
$ws.Range("S4").Value = 0.07808271846324491
$ws.Range("T4").Value = 0.09483404869534075

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Cadm3"
$ws.Range("C5").Value = "Cadm1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.199582666666667
$ws.Range("H5").Value = 12.598748
$ws.Range("I5").Value = 0.3470882463138872
$ws.Range("J5").Value = 0.3745187327618799
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.589755
$ws.Range("N5").Value = 19.17951
$ws.Range("O5").Value = 0.3347035199121805
$ws.Range("P5").Value = 0.2511566892559376
$ws.Range("Q5").Value = 40.27296887558001
$ws.Range("R5").Value = 241.63781325348
$ws.Range("S5").Value = 0.116171657761404
$ws.Range("T5").Value = 0.09406288498480302

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cadm3"
$ws.Range("C6").Value = "Cadm1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 4.362043666666667
$ws.Range("H6").Value = 13.086131
$ws.Range("I6").Value = 0.3605153670685209
$ws.Range("J6").Value = 0.3890070028288486
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.515984
$ws.Range("N6").Value = 7.547952
$ws.Range("O6").Value = 0.08781336966822693
$ws.Range("P6").Value = 0.09884082726736673
$ws.Range("Q6").Value = 10.97483207263467
$ws.Range("R6").Value = 98.77348865371199
$ws.Range("S6").Value = 0.03165806919946455
$ws.Range("T6").Value = 0.03844977397240227

# Row 7: FAPs -> M1
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cadm3"
$ws.Range("C7").Value = "Cadm1"
$ws.Range("D7").Value = "M1"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4.362043666666667
$ws.Range("H7").Value = 13.086131
$ws.Range("I7").Value = 0.3605153670685209
$ws.Range("J7").Value = 0.3890070028288486
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.10016866666667
$ws.Range("N7").Value = 30.300506
$ws.Range("O7").Value = 0.3525180783492434
$ws.Range("P7").Value = 0.3967867150797739
$ws.Range("Q7").Value = 44.05737676469845
$ws.Range("R7").Value = 396.516390882286
$ws.Range("S7").Value = 0.1270881844143671
$ws.Range("T7").Value = 0.1543528107954872

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Cadm3"
$ws.Range("C8").Value = "Cadm1"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.362043666666667
$ws.Range("H8").Value = 13.086131
$ws.Range("I8").Value = 0.3605153670685209
$ws.Range("J8").Value = 0.3890070028288486
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 6.445583666666667
$ws.Range("N8").Value = 19.336751
$ws.Range("O8").Value = 0.2249650320703493
$ws.Range("P8").Value = 0.2532157683969216
$ws.Range("Q8").Value = 28.11591741115344
$ws.Range("R8").Value = 253.043256700381
$ws.Range("S8").Value = 0.08110335111442354
$ws.Range("T8").Value = 0.09850270713309037

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Cadm3"
$ws.Range("C9").Value = "Cadm1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.362043666666667
$ws.Range("H9").Value = 13.086131
$ws.Range("I9").Value = 0.3605153670685209
$ws.Range("J9").Value = 0.3890070028288486
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.589755
$ws.Range("N9").Value = 19.17951
$ws.Range("O9").Value = 0.3347035199121805
$ws.Range("P9").Value = 0.2511566892559376
$ws.Range("Q9").Value = 41.830930062635
$ws.Range("R9").Value = 250.98558037581
$ws.Range("S9").Value = 0.1206657623402658
$ws.Range("T9").Value = 0.09770171092786879

# Row 10: M1 -> ECs
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Cadm3"
$ws.Range("C10").Value = "Cadm1"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7916129999999999
$ws.Range("H10").Value = 2.374839
$ws.Range("I10").Value = 0.0654254457496749
$ws.Range("J10").Value = 0.07059603801849913
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.515984
$ws.Range("N10").Value = 7.547952
$ws.Range("O10").Value = 0.08781336966822693
$ws.Range("P10").Value = 0.09884082726736673
$ws.Range("Q10").Value = 1.991685642192
$ws.Range("R10").Value = 17.925170779728
$ws.Range("S10").Value = 0.005745228853324729
$ws.Range("T10").Value = 0.006977770799546927

# Row 11: M1 -> M1
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = "Cadm3"
$ws.Range("C11").Value = "Cadm1"
$ws.Range("D11").Value = "M1"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7916129999999999
$ws.Range("H11").Value = 2.374839
$ws.Range("I11").Value = 0.0654254457496749
$ws.Range("J11").Value = 0.07059603801849913
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 10.10016866666667
$ws.Range("N11").Value = 30.300506
$ws.Range("O11").Value = 0.3525180783492434
$ws.Range("P11").Value = 0.3967867150797739
$ws.Range("Q11").Value = 7.995424818725999
$ws.Range("R11").Value = 71.95882336853398
$ws.Range("S11").Value = 0.02306365241081807
$ws.Range("T11").Value = 0.0280115700230071

# Row 12: M1 -> M2
$ws.Range("A12").Value = "M1"
$ws.Range("B12").Value = "Cadm3"
$ws.Range("C12").Value = "Cadm1"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7916129999999999
$ws.Range("H12").Value = 2.374839
$ws.Range("I12").Value = 0.0654254457496749
$ws.Range("J12").Value = 0.07059603801849913
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 6.445583666666667
$ws.Range("N12").Value = 19.336751
$ws.Range("O12").Value = 0.2249650320703493
$ws.Range("P12").Value = 0.2532157683969216
$ws.Range("Q12").Value = 5.102407823120999
$ws.Range("R12").Value = 45.92167040808899
$ws.Range("S12").Value = 0.01471843750129251
$ws.Range("T12").Value = 0.01787603001263255

# Row 13: M1 -> sCs
$ws.Range("A13").Value = "M1"
$ws.Range("B13").Value = "Cadm3"
$ws.Range("C13").Value = "Cadm1"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7916129999999999
$ws.Range("H13").Value = 2.374839
$ws.Range("I13").Value = 0.0654254457496749
$ws.Range("J13").Value = 0.07059603801849913
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.589755
$ws.Range("N13").Value = 19.17951
$ws.Range("O13").Value = 0.3347035199121805
$ws.Range("P13").Value = 0.2511566892559376
$ws.Range("Q13").Value = 7.591374724814999
$ws.Range("R13").Value = 45.54824834889
$ws.Range("S13").Value = 0.0218981269842396
$ws.Range("T13").Value = 0.01773066718331255

# Row 14: M2 -> ECs
$ws.Range("A14").Value = "M2"
$ws.Range("B14").Value = "Cadm3"
$ws.Range("C14").Value = "Cadm1"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.08766133333333333
$ws.Range("H14").Value = 0.262984
$ws.Range("I14").Value = 0.007245057633394309
$ws.Range("J14").Value = 0.007817636674425921
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.515984
$ws.Range("N14").Value = 7.547952
$ws.Range("O14").Value = 0.08781336966822693
$ws.Range("P14").Value = 0.09884082726736673
$ws.Range("Q14").Value = 0.2205545120853333
$ws.Range("R14").Value = 1.984990608768
$ws.Range("S14").Value = 0.0006362129242288638
$ws.Range("T14").Value = 0.0007727016761759637

# Row 15: M2 -> M1
$ws.Range("A15").Value = "M2"
$ws.Range("B15").Value = "Cadm3"
$ws.Range("C15").Value = "Cadm1"
$ws.Range("D15").Value = "M1"
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.08766133333333333
$ws.Range("H15").Value = 0.262984
$ws.Range("I15").Value = 0.007245057633394309
$ws.Range("J15").Value = 0.007817636674425921
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 10.10016866666667
$ws.Range("N15").Value = 30.300506
$ws.Range("O15").Value = 0.3525180783492434
$ws.Range("P15").Value = 0.3967867150797739
$ws.Range("Q15").Value = 0.8853942522115555
$ws.Range("R15").Value = 7.968548269904
$ws.Range("S15").Value = 0.002554013794453679
$ws.Range("T15").Value = 0.003101934375732629

# Row 16: M2 -> M2
$ws.Range("A16").Value = "M2"
$ws.Range("B16").Value = "Cadm3"
$ws.Range("C16").Value = "Cadm1"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.08766133333333333
$ws.Range("H16").Value = 0.262984
$ws.Range("I16").Value = 0.007245057633394309
$ws.Range("J16").Value = 0.007817636674425921
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 6.445583666666667
$ws.Range("N16").Value = 19.336751
$ws.Range("O16").Value = 0.2249650320703493
$ws.Range("P16").Value = 0.2532157683969216
$ws.Range("Q16").Value = 0.5650284583315556
$ws.Range("R16").Value = 5.085256124983999
$ws.Range("S16").Value = 0.00162988462284808
$ws.Range("T16").Value = 0.001979548877562715

# Row 17: M2 -> sCs
$ws.Range("A17").Value = "M2"
$ws.Range("B17").Value = "Cadm3"
$ws.Range("C17").Value = "Cadm1"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.08766133333333333
$ws.Range("H17").Value = 0.262984
$ws.Range("I17").Value = 0.007245057633394309
$ws.Range("J17").Value = 0.007817636674425921
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 9.589755
$ws.Range("N17").Value = 19.17951
$ws.Range("O17").Value = 0.3347035199121805
$ws.Range("P17").Value = 0.2511566892559376
$ws.Range("Q17").Value = 0.84065070964
$ws.Range("R17").Value = 5.04390425784
$ws.Range("S17").Value = 0.002424946291863688
$ws.Range("T17").Value = 0.001963451744954613

# Row 18: sCs -> ECs
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Cadm3"
$ws.Range("C18").Value = "Cadm1"
$ws.Range("D18").Value = "ECs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 2.658566
$ws.Range("H18").Value = 5.317132
$ws.Range("I18").Value = 0.2197258832345227
$ws.Range("J18").Value = 0.1580605897163464
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 2.515984
$ws.Range("N18").Value = 7.547952
$ws.Range("O18").Value = 0.08781336966822693
$ws.Range("P18").Value = 0.09884082726736673
$ws.Range("Q18").Value = 6.688909518944
$ws.Range("R18").Value = 40.133457113664
$ws.Range("S18").Value = 0.01929487021015081
$ws.Range("T18").Value = 0.01562283944593152

# Row 19: sCs -> M1
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Cadm3"
$ws.Range("C19").Value = "Cadm1"
$ws.Range("D19").Value = "M1"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 2.658566
$ws.Range("H19").Value = 5.317132
$ws.Range("I19").Value = 0.2197258832345227
$ws.Range("J19").Value = 0.1580605897163464
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 10.10016866666667
$ws.Range("N19").Value = 30.300506
$ws.Range("O19").Value = 0.3525180783492434
$ws.Range("P19").Value = 0.3967867150797739
$ws.Range("Q19").Value = 26.85196501146534
$ws.Range("R19").Value = 161.111790068792
$ws.Range("S19").Value = 0.07745734612142419
$ws.Range("T19").Value = 0.06271634217712098

# Row 20: sCs -> M2
$ws.Range("A20").Value = "sCs"
$ws.Range("B20").Value = "Cadm3"
$ws.Range("C20").Value = "Cadm1"
$ws.Range("D20").Value = "M2"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 2.658566
$ws.Range("H20").Value = 5.317132
$ws.Range("I20").Value = 0.2197258832345227
$ws.Range("J20").Value = 0.1580605897163464
$ws.Range("K20").Value = 3
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 6.445583666666667
$ws.Range("N20").Value = 19.336751
$ws.Range("O20").Value = 0.2249650320703493
$ws.Range("P20").Value = 0.2532157683969216
$ws.Range("Q20").Value = 17.13600958635533
$ws.Range("R20").Value = 102.816057518132
$ws.Range("S20").Value = 0.04943064036854022
$ws.Range("T20").Value = 0.04002343367829522

# Row 21: sCs -> sCs
$ws.Range("A21").Value = "sCs"
$ws.Range("B21").Value = "Cadm3"
$ws.Range("C21").Value = "Cadm1"
$ws.Range("D21").Value = "sCs"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 2.658566
$ws.Range("H21").Value = 5.317132
$ws.Range("I21").Value = 0.2197258832345227
$ws.Range("J21").Value = 0.1580605897163464
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 9.589755
$ws.Range("N21").Value = 19.17951
$ws.Range("O21").Value = 0.3347035199121805
$ws.Range("P21").Value = 0.2511566892559376
$ws.Range("Q21").Value = 25.49499659133
$ws.Range("R21").Value = 101.97998636532
$ws.Range("S21").Value = 0.07354302653440753
$ws.Range("T21").Value = 0.03969797441499866

